# Generate Report for Handoff
# Adds two new handed-off files (b470feb8... and df4d5edc...) as new rows
# to the Overview sheet and the zh-cn / de-de per-language sheets, expanding
# the backing tables accordingly.

$wb = $excel.ActiveWorkbook

$ovName    = "Overview"
$zhName    = "zh-cn"
$deName    = "de-de"

$wsOv = $wb.Worksheets.Item($ovName)
$wsZh = $wb.Worksheets.Item($zhName)
$wsDe = $wb.Worksheets.Item($deName)

# ---------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel coercing
# look-alike numbers/dates/booleans - prefix with an apostrophe (which
# Excel strips) to force the "text" interpretation, then reset the style
# back to Normal/"@" free defaults.
# ---------------------------------------------------------------------

function Set-Text($rng, $val) {
    if ($val -eq "") {
        $rng.Value = ""
    } else {
        $rng.Value = "'" + $val
    }
}

# ======================= Overview sheet =======================
$tblOv = $wsOv.ListObjects.Item(1)
$tblOv.ListRows.Add() | Out-Null
$tblOv.ListRows.Add() | Out-Null

$ovRows = @(
    @{ Row=6; Id="b470feb8-1f96-4f95-a14c-a56266fcd51a"; Status="Ready for handoff"; Date="2016-08-18 04:40:34" },
    @{ Row=7; Id="df4d5edc-1fd0-4975-ac1d-354b2bb5fc19"; Status="Ready for handoff"; Date="2016-08-18 04:40:34" }
)

foreach ($r in $ovRows) {
    $row = $r.Row
    $fname = $r.Id + ".md"
    $path  = "e2e\" + $fname

    Set-Text $wsOv.Range("A$row") $fname
    Set-Text $wsOv.Range("C$row") ".md"
    Set-Text $wsOv.Range("D$row") ""
    Set-Text $wsOv.Range("E$row") $r.Status
    Set-Text $wsOv.Range("F$row") $r.Status
    Set-Text $wsOv.Range("G$row") $r.Date

    # B column carries the hyperlink (display text = path, target = github URL)
    $wsOv.Range("B$row").Value = "'" + $path
    $wsOv.Hyperlinks.Add($wsOv.Range("B$row"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/$fname", "", "", $path) | Out-Null

    $wsOv.Range("A$row").Style = "Normal"
    $wsOv.Range("B$row").Style = "HyperLink"
    $wsOv.Range("C$row").Style = "Normal"
    $wsOv.Range("E$row").Style = "Normal"
    $wsOv.Range("F$row").Style = "Normal"
    $wsOv.Range("G$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ======================= Per-language sheets =======================
function Fill-LangSheet($ws, $tbl, $langSuffix, $hoHash1, $hoHash2, $dates) {
    $tbl.ListRows.Add() | Out-Null
    $tbl.ListRows.Add() | Out-Null

    $entries = @(
        @{ Row=6; Id="b470feb8-1f96-4f95-a14c-a56266fcd51a"; Hash=$hoHash1; Date=$dates[0] },
        @{ Row=7; Id="df4d5edc-1fd0-4975-ac1d-354b2bb5fc19"; Hash=$hoHash2; Date=$dates[1] }
    )

    foreach ($e in $entries) {
        $row = $e.Row
        $fname = $e.Id + ".md"
        $xlf   = $e.Id + "." + $e.Hash + "." + $langSuffix + ".xlf"

        $ws.Range("A$row").Value = "'" + $fname
        Set-Text $ws.Range("B$row") ".md"
        Set-Text $ws.Range("C$row") "Ready for handoff"
        Set-Text $ws.Range("D$row") "e2e"
        Set-Text $ws.Range("E$row") "ht"
        Set-Text $ws.Range("F$row") "False"
        Set-Text $ws.Range("G$row") $xlf
        Set-Text $ws.Range("H$row") $e.Date
        Set-Text $ws.Range("I$row") ""
        Set-Text $ws.Range("J$row") ""
        Set-Text $ws.Range("K$row") "0001-01-01 00:00:00"
        Set-Text $ws.Range("L$row") ""
        Set-Text $ws.Range("M$row") "True"
        Set-Text $ws.Range("N$row") ""
        Set-Text $ws.Range("O$row") "False"
        Set-Text $ws.Range("P$row") ""

        $ws.Hyperlinks.Add($ws.Range("A$row"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/$fname", "", "", $fname) | Out-Null

        $ws.Range("A$row").Style = "HyperLink"
        foreach ($col in @("B","C","D","E","F","G","I","J","L","M","N","O","P")) {
            $ws.Range("$col$row").Style = "Normal"
        }
        $ws.Range("H$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"
        $ws.Range("K$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    }
}

$tblZh = $wsZh.ListObjects.Item(1)
Fill-LangSheet $wsZh $tblZh "zh-cn" "eeb70b3e6879e58c88c540615cd3c66d88d6f7a9" "8868acd87d91866d9248d4819a9c9a3445087d29" @("2016-08-18 04:40:29","2016-08-18 04:40:29")

$tblDe = $wsDe.ListObjects.Item(1)
Fill-LangSheet $wsDe $tblDe "de-de" "eeb70b3e6879e58c88c540615cd3c66d88d6f7a9" "8868acd87d91866d9248d4819a9c9a3445087d29" @("2016-08-18 04:40:34","2016-08-18 04:40:34")

Write-Output "Handoff report rows appended."
